# apuesta_dia.xlsx — replace the "bet of the day" table with the next day's
# picks: new date, four new DET/BOS/NOP/MIN matchups (two bets each) instead
# of the previous three ORL/PHI/CHI matchups (three bets each), new running
# totals/final odds, and the new Stake bet-slip link.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header date moves from 2025-04-09 to 2025-04-11 ---
$ws.Range("A1").Value = 45758

# --- The three sub-total rows used to sit at rows 5, 9, 13 (style-wise the
#     "match" column (A) and "info" column (D) are blank on a sub-total row).
#     In the new layout there's one more sub-total row and they land on rows
#     4, 7, 10, 13, so A4/D4, A7/D7, A10/D10 (which held old match-name /
#     info text) must be cleared out. ---
$ws.Range("A4").ClearContents()
$ws.Range("D4").ClearContents()
$ws.Range("A7").ClearContents()
$ws.Range("D7").ClearContents()
$ws.Range("A10").ClearContents()
$ws.Range("D10").ClearContents()

# --- Row heights: rows 3 & 4 shrink to the regular 24.95pt data-row height,
#     row 7 grows to the 33.75pt sub-total-row height (it's now a "CUOTA
#     TOTAL =" row instead of a pick row). ---
$ws.Rows.Item(3).RowHeight = 24.95
$ws.Rows.Item(4).RowHeight = 33.75
$ws.Rows.Item(7).RowHeight = 33.75

# --- New table content ---
# Match 1: DET vs MIL
$ws.Range("A2").Value = "DET vs MIL"
$ws.Range("B2").Value = "Giannis Antetokounmpo - Sobre 8.5 Goles de campo realizados"
$ws.Range("C2").Value = "1,53x"
$ws.Range("D2").Value = "Cumplio en 10/10 (pueden bajar la linea a 7,5) https://prnt.sc/4eC8v1GOiTYK"

$ws.Range("A3").Value = "DET vs MIL"
$ws.Range("B3").Value = "AJ Green - Sobre 3.5 Gol de Campo intentado"
$ws.Range("C3").Value = "1,37x"
$ws.Range("D3").Value = "Cumplio en 10/10 https://prnt.sc/EzBJAad6nWHG"

$ws.Range("B4").Value = "CUOTA TOTAL ="
$ws.Range("C4").Value = "2,26x"

# Match 2: BOS vs CHA
$ws.Range("A5").Value = "BOS vs CHA"
$ws.Range("B5").Value = "Jrue Holiday - Sobre 3.5 Dobles intentados"
$ws.Range("C5").Value = "1,59x"
$ws.Range("D5").Value = "Cumplio en 9/10 (hizo 2) https://prnt.sc/O_qWGBROzl5t"

$ws.Range("A6").Value = "BOS vs CHA"
$ws.Range("B6").Value = "Kristaps Porzingis - Sobre 6.5 Dobles intentados"
$ws.Range("C6").Value = "1,55x"
$ws.Range("D6").Value = "Cumplio en 9/10 (hizo 6, pueden bajar a 5,5 intentos) https://prnt.sc/nj7_Gbfhl2WG"

$ws.Range("B7").Value = "CUOTA TOTAL ="
$ws.Range("C7").Value = "2,65x"

# Match 3: NOP vs MIA
$ws.Range("A8").Value = "NOP vs MIA"
$ws.Range("B8").Value = "Davion Mitchell - Sobre 9.5 Puntos"
$ws.Range("C8").Value = "1,48x"
$ws.Range("D8").Value = "Cumplio en 9/10 https://prnt.sc/T7LPp3-rSUcj"

$ws.Range("A9").Value = "NOP vs MIA"
$ws.Range("B9").Value = "Andrew Wiggins - Sobre 0.5 Faltas personales"
$ws.Range("C9").Value = "1,33x"
$ws.Range("D9").Value = "Cumplio en 10/10 https://prnt.sc/0smaNnDUB8xN"

$ws.Range("B10").Value = "CUOTA TOTAL ="
$ws.Range("C10").Value = "2,08x"

# Match 4: MIN vs BKN
$ws.Range("A11").Value = "MIN vs BKN"
$ws.Range("B11").Value = "Nickeil Alexander-Walker - Sobre 0.5 Tiros Libres intentados"
$ws.Range("C11").Value = "1,83x"
$ws.Range("D11").Value = "Cumplio en 9/10 https://prnt.sc/ycK4i6bfu2W2"

$ws.Range("A12").Value = "MIN vs BKN"
$ws.Range("B12").Value = "Rudy Gobert - Sobre 1.5 Faltas personales"
$ws.Range("C12").Value = "1,67x"
$ws.Range("D12").Value = "Cumplio en 9/10 https://prnt.sc/sDnjSUIsjhkH"

$ws.Range("B13").Value = "CUOTA TOTAL ="
$ws.Range("C13").Value = "3,35x"

# Grand total row
$ws.Range("B14").Value = "CUOTA APUESTA TOTAL ="
$ws.Range("C14").Value = "53,02x"
$ws.Range("D14").Value = "https://stake.com/sports/home?betId=03756cc3-756b-426c-8d91-1f39b94d0e0a&modal=bet"

# --- The two trailing blank rows (15, 16) left over from the old, taller
#     table are removed entirely. ---
$ws.Rows.Item(15).Delete()
$ws.Rows.Item(15).Delete()

# --- Selection moves off the (now sub-total) B4 cell to D17, matching the
#     author's last on-screen click before saving. ---
$ws.Range("D17").Select()
